# Auto commit at 2025-10-09 8:30:06.79
# Updates the Metrics source values, which ripple via formulas into the
# "today" sheet, and moves the active-sheet/selection state from
# "today"!G7 to "Metrics"!D15 (today's selection moves on to G11).

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# --- Update the raw metric values on the Metrics sheet (B2:B13) ---
$metrics.Range("B2").Value  = 109245.07
$metrics.Range("B3").Value  = 92517.37000000001
$metrics.Range("B4").Value  = 41930.730000000003
$metrics.Range("B5").Value  = 4218
$metrics.Range("B6").Value  = 4476376.540000001
$metrics.Range("B7").Value  = 3782336.0399999996
$metrics.Range("B8").Value  = 1312532.8700000001
$metrics.Range("B9").Value  = 173219
$metrics.Range("B10").Value = 32941700.340999827
$metrics.Range("B11").Value = 31057557.560000006
$metrics.Range("B12").Value = 11594241.76
$metrics.Range("B13").Value = 1270846

# --- Move the selection on "today" off of G7 and onto G11 ---
$today.Activate()
$today.Range("G11").Select()

# --- Make "Metrics" the active sheet/tab, selecting D15 ---
$metrics.Activate()
$metrics.Range("D15").Select()
